$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(
        "Golang Developer - Onsite Role",
        "https://www.dice.com/job-detail/2fca2109-3102-4424-bcfa-92d1fbb17518?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang",
        "Dallas, Texas",
        "Contract, Third Party",
        "Depends on Experience",
        "Navtech, LLC"
    ),
    @(
        "Sr.Software Developer(Golang, SQL, Postgres)",
        "https://www.dice.com/job-detail/860a5fe9-3cd6-43de-a64a-b4435aa147e5?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang",
        "Remote",
        "Contract",
        "Depends on Experience",
        "Nobl Q, LLC."
    ),
    @(
        "Golang Engineer (AWS & Java) - Plano, TX / McLean, VA / Richmond, VA (3 days onsite)",
        "https://www.dice.com/job-detail/eb6d162f-dc2a-4432-9a89-9e50a884c75f?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang",
        "Hybrid in Plano, Texas",
        "Contract, Third Party",
        "Depends on Experience",
        "RKube Inc"
    ),
    @(
        "Golang Developer",
        "https://www.dice.com/job-detail/ad009a2b-e4d0-4713-a240-7edb44f0a065?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang",
        "McLean, Virginia",
        "Contract",
        "$50 - $60",
        "S Linx LLC"
    ),
    @(
        "Java Backend Developer",
        "https://www.dice.com/job-detail/9c157b82-9fe0-4fec-a6e4-f5da6d9fec4c?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang",
        "Hybrid in McLean, Virginia",
        "Contract, Third Party",
        "50 - 60",
        "Gtech LLC"
    )
)

$startRow = 132
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
